$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the CaseFilesTab query (cell B4) to split the generic .gz rule into
# specific .fastq.gz / .vcf.gz / .tsv.gz rules.
$newQuery = @'
SELECT 
    DISTINCT cf.file_name AS "File Name",
    CASE
        WHEN cf.file_name LIKE '%.bai' THEN 'bai'
        WHEN cf.file_name LIKE '%.bam' THEN 'bam'
        WHEN cf.file_name LIKE '%.csv' THEN 'csv'
        WHEN cf.file_name LIKE '%.doc' THEN 'doc'
        WHEN cf.file_name LIKE '%.docx' THEN 'docx'
           WHEN cf.file_name LIKE '%.fastq.gz' THEN 'fastq'
        WHEN cf.file_name LIKE '%.vcf.gz' THEN 'vcf'
        WHEN cf.file_name LIKE '%.tsv.gz' THEN 'tsv'
        WHEN cf.file_name LIKE '%.pdf' THEN 'pdf'
        WHEN cf.file_name LIKE '%.rtf' THEN 'rtf'
        WHEN cf.file_name LIKE '%.tbi' THEN 'tbi'
        WHEN cf.file_name LIKE '%.tif' THEN 'tif'
        WHEN cf.file_name LIKE '%.xls' THEN 'xls'
        WHEN cf.file_name LIKE '%.xlsx' THEN 'xlsx'
        ELSE 'Unknown'
    END AS "Format",
    cf.file_type AS "File Type",
    CASE     
    WHEN cf.file_size >= 1024 * 1024 * 1024 THEN 
        ROUND(cf.file_size / (1024.0 * 1024.0 * 1024.0), 2) || ' GB'
    WHEN cf.file_size >= 1024 * 1024 THEN 
        ROUND(cf.file_size / (1024.0 * 1024.0), 2) || ' MB'
    WHEN cf.file_size >= 1024 THEN 
        ROUND(cf.file_size / 1024.0, 2) || ' KB'
    ELSE 
        cf.file_size || ' Bytes'
    END AS "Size",
    CASE 
        WHEN smp.sample_id IS NOT NULL THEN 'sample'
        WHEN c.case_record_id IS NOT NULL THEN 'case'
        ELSE 'Unknown'
    END AS "Association",
    cf.file_description AS "Description",
    smp.sample_id AS "Sample ID",
    c.case_record_id AS "Case ID",
    dmg.breed AS "Breed",
    diag.disease_term AS "Diagnosis"
FROM 
    df_case_file cf
LEFT JOIN 
    df_sample smp ON cf."sample.sample_id" = smp.sample_id
LEFT JOIN 
    df_case c ON (smp."case.case_record_id" = c.case_record_id OR cf."case.case_record_id" = c.case_record_id)
LEFT JOIN 
    df_study st ON c."study.clinical_study_designation" = st.clinical_study_designation
LEFT JOIN 
    df_demographic dmg ON dmg."case.case_record_id" = c.case_record_id
LEFT JOIN 
    df_diagnosis diag ON diag."case.case_record_id" = c.case_record_id
WHERE
    st.clinical_study_designation = 'GLIOMA01'
    AND cf."case.case_record_id" IS NOT NULL 
ORDER BY 
    cf.file_name ASC
LIMIT 100;
'@

$ws.Range("B4").Value = $newQuery

# Setting the longer, multi-line query text causes the wrap-text row to
# auto-grow; restore the row's original (already-maxed-out) height so the
# row geometry is left exactly as it was.
$ws.Rows.Item(4).RowHeight = 409.6

# Move the active selection from B5 to B4, matching the updated sheetView.
$ws.Range("B4").Select()
